$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.622815847396851
$ws.Range("B1").Value = 1.616479992866516
$ws.Range("C1").Value = 1.948170185089111
$ws.Range("D1").Value = 3.683640718460083
$ws.Range("E1").Value = 4.388773441314697
